$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 6 (shifts the existing 2025-Jun..2018-Ene data rows
# and the footer rows down by one), matching table row B5:G95 -> B5:G96.
$ws.Rows.Item(6).Insert(-4121)

# Pick up the banded-row formatting (fill/border/number-format) from the row
# that is now two rows below (originally row 7, the "May." banded row) so the
# newly inserted row looks like a normal table data row instead of inheriting
# the header's format.
$ws.Range("B8:G8").Copy()
$ws.Range("B6:G6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new July 2025 data row.
$ws.Range("B6").Value2 = 2025
$ws.Range("C6").Value2 = "Jul."
$ws.Range("D6").Value2 = 30935
$ws.Range("E6").Value2 = 4284
$ws.Range("F6").Value2 = 21100
$ws.Range("G6").Value2 = 444

# The June 2025 row (now row 7) had its "Carga" figure revised.
$ws.Range("D7").Value2 = 24801

# Grow the table/autofilter range so it includes the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B5:G96"))

# Bump the "last updated" caption.
$ws.Range("B97").Value2 = "Actualización: Julio 2025."
